$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Fecha de actualizacion" day: 03 -> 07. Only the trailing "3" is
#    touched (replaced by a new run holding "7") so the rest of the date's
#    pre-existing run-split structure is left untouched.
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "03"
$find.Execute() | Out-Null
$dayRange = $d.Range($find.Parent.Start, $find.Parent.End)
$tail = $d.Range($dayRange.Start + 1, $dayRange.End)

$sevenRunXml = @"
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r>
  <w:rPr>
    <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
    <w:sz w:val="22"/>
    <w:szCs w:val="22"/>
    <w:lang w:val="es-CO"/>
  </w:rPr>
  <w:t>7</w:t>
</w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$tail.InsertXML($sevenRunXml)

# ---------------------------------------------------------------------------
# 2) Intro paragraph: table suffix "raw" -> "transformed"
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("consolidado_junio2023_mayo2024_raw", $false, $false, $false, $false, $false, $true, 1, $false, "consolidado_junio2023_mayo2024_transformed", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Collapse the three "relacion entre el tipo de miembro y ..." bullet
#    questions (distancias / estaciones / rutas) down to a single empty
#    list-style paragraph (keeps pStyle, drops numPr).
# ---------------------------------------------------------------------------
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*distancias recorridas en los viajes?*") { $startPara = $d.Paragraphs.Item($i) }
    if ($t -like "*alguna relaci*n entre el tipo de miembro y las rutas utilizadas?*") { $endPara = $d.Paragraphs.Item($i) }
}

$rng = $d.Range($startPara.Range.Start, $endPara.Range.End)

$emptyListXml = @"
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Prrafodelista"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
      <w:lang w:val="es-CO"/>
    </w:rPr>
  </w:pPr>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$rng.InsertXML($emptyListXml)

# ---------------------------------------------------------------------------
# 4) Merge the now-empty bold paragraph with the "Analisis descriptivo"
#    heading paragraph right after it (delete the paragraph mark between
#    them so the heading run lives inside the first paragraph's pPr).
# ---------------------------------------------------------------------------
$headingPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)
    if ($t -eq "Análisis descriptivo") { $headingPara = $d.Paragraphs.Item($i); break }
}
$prevPara = $headingPara.Previous()
$prevText = $prevPara.Range.Text.TrimEnd([char]13)
if ($prevText -eq "") {
    $mark = $d.Range($prevPara.Range.End - 1, $prevPara.Range.End)
    $mark.Delete(1, 1) | Out-Null
}

$d.Save()
